# The "Coral (hard and soft)" row (row 3) is removed from the taxonomic
# richness table; everything below it shifts up by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(3).Delete() | Out-Null

# The AutoFilter range shrinks along with the table (was A9:C18, now A8:C17).
$wb.Names("_xlnm._FilterDatabase").RefersTo = "=Sheet1!`$A`$8:`$C`$17"

# Restore the selection left behind by the author's last edit.
$ws.Range("G11").Select() | Out-Null
